$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark from the last (empty) paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Add first-line indent (720 twips = 36pt) to the first paragraph.
$p1 = $d.Paragraphs.Item(1)
$p1.Format.FirstLineIndent = 36

# 3. Re-insert the "_GoBack" bookmark collapsed at the very start of the
#    first paragraph (before its run), matching the target layout.
#    A zero-length range placed exactly at document position 0 gets
#    mis-anchored by the engine, landing its end in the following
#    paragraph, so work around it: insert a throw-away character at the
#    very start, anchor the bookmark right after it (a safe, unambiguous
#    position within paragraph 1), then remove the throw-away character.
$tmp = $d.Range(0, 0)
$tmp.InsertBefore("Z")
$anchor = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $anchor)
$d.Range(0, 1).Delete()
